# More work on excel export. Added coloration of section header.
#
# This adds a 4th line ("Modeled the skill") to the first skill-building
# section, and adds a small column of annotation text (in column H)
# explaining the layout of each section/header to whoever is looking
# at the template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row right after row 3 so the first section gains a 4th
# list item; this pushes the two sections below it down by one row,
# which Excel keeps in sync with their existing merged header cells.
$ws.Rows("4").Insert()
$ws.Range("A4").Value = "Modeled the skill"

# Annotation notes explaining the template's layout, placed to the
# right of the form so they don't interfere with the printable area.
$ws.Range("H1").Value = "<- Each section has it's own header"
$ws.Range("H2").Value = "<- the skill is a larger font"
$ws.Range("H3").Value = "<- each selection is then listed under the skill"

# Leave the active selection where the author left off editing.
$ws.Range("H4").Select() | Out-Null
